# Add the new I0 ("I0") and IF ("IF") columns to the king_john sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns (row 1), matching the style used by the
# other header cells (B1:H1). Copy/PasteSpecial(Formats) reuses the existing
# header cell style instead of fabricating a brand-new one.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for rows 2-33 for columns I (I0) and J (IF).
$i0Values = @(6, 5, 6, 7, 5, 7, 8, 6, 5, 7, 10, 7, 7, 6, 7, 7, 4, 7, 6, 7, 7, 10, 6, 8, 6, 6, 6, 5, 5, 5, 4, 4)
$ifValues = @(7, 8, 9, 9, 7, 8, 8, 7, 7, 8, 10, 8, 8, 7, 7, 8, 6, 7, 7, 9, 7, 10, 7, 8, 7, 6, 6, 7, 7, 6, 5, 4)

for ($r = 2; $r -le 33; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($r, 10).Value = $ifValues[$idx]
}
